$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12883.25
$ws.Range("I6").Value = 16843.666
$ws.Range("J6").Value = 1002
$ws.Range("K6").Value = 50530.99800000001
$ws.Range("L6").Value = 3006
$ws.Range("M6").Value = -50418.99800000001
$ws.Range("N6").Value = -3230
$ws.Range("H33").Value = 391.10715
$ws.Range("I33").Value = 353.30435
$ws.Range("K33").Value = 353.30435
$ws.Range("M33").Value = -124.30435
$ws.Range("H64").Value = 4087.5
$ws.Range("J64").Value = 3958.5715
$ws.Range("L64").Value = 3958.5715
$ws.Range("N64").Value = -4454.5715
$ws.Range("H67").Value = 4087.5
$ws.Range("J67").Value = 3958.5715
$ws.Range("L67").Value = 3958.5715
$ws.Range("N67").Value = -5674.5715
$ws.Range("H106").Value = 9646.467000000001
$ws.Range("I106").Value = 10592.077
$ws.Range("K106").Value = 10592.077
$ws.Range("M106").Value = -9961.076999999999
$ws.Range("H113").Value = 2617.375
$ws.Range("J113").Value = 2664.75
$ws.Range("L113").Value = 2664.75
$ws.Range("N113").Value = -9172.75
$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3861.87
$ws.Range("I32").Value = 3724.8667
$ws.Range("K32").Value = 3724.8667
$ws.Range("M32").Value = -3437.8667
$ws.Range("H61").Value = 52632840
$ws.Range("I61").Value = 76924040
$ws.Range("J61").Value = 1906.6666
$ws.Range("K61").Value = 76924040
$ws.Range("L61").Value = 1906.6666
$ws.Range("M61").Value = -76923828
$ws.Range("N61").Value = -2330.6666
$ws.Range("H74").Value = 1401.8823
$ws.Range("I74").Value = 881.5714
$ws.Range("K74").Value = 881.5714
$ws.Range("M74").Value = -7.57140000000004
$ws.Range("H77").Value = 1401.8823
$ws.Range("I77").Value = 881.5714
$ws.Range("K77").Value = 4407.857
$ws.Range("M77").Value = -39.85699999999997
$ws.Range("H97").Value = 718.55554
$ws.Range("I97").Value = 658.375
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 658.375
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -162.375
$ws.Range("N97").Value = -2192
$ws.Range("H132").Value = 1897.8292
$ws.Range("I132").Value = 1322.0333
$ws.Range("K132").Value = 3966.0999
$ws.Range("M132").Value = -1436.0999
$ws.Range("H136").Value = 52632840
$ws.Range("I136").Value = 76924040
$ws.Range("J136").Value = 1906.6666
$ws.Range("K136").Value = 230772120
$ws.Range("L136").Value = 5719.9998
$ws.Range("M136").Value = -230769570
$ws.Range("N136").Value = -10819.9998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 24000
$ws.Range("I93").Value = 18000
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 18000
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -16128
$ws.Range("N93").Value = -33744
$ws.Range("H99").Value = 37038144
$ws.Range("I99").Value = 45455628
$ws.Range("J99").Value = 1212.2
$ws.Range("K99").Value = 45455628
$ws.Range("L99").Value = 1212.2
$ws.Range("M99").Value = -45454130
$ws.Range("N99").Value = -4208.2
$ws.Range("H134").Value = 1371.92
$ws.Range("I134").Value = 877.4375
$ws.Range("K134").Value = 2632.3125
$ws.Range("M134").Value = -97.3125

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 166668220
$ws.Range("I16").Value = 166668220
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 166668220
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -166667933
$ws.Range("H31").Value = 1695.96
$ws.Range("I31").Value = 1573.2354
$ws.Range("K31").Value = 1573.2354
$ws.Range("M31").Value = -1278.2354
$ws.Range("H34").Value = 1695.96
$ws.Range("I34").Value = 1573.2354
$ws.Range("K34").Value = 1573.2354
$ws.Range("M34").Value = -1371.2354
$ws.Range("H86").Value = 4483726
$ws.Range("I86").Value = 11153860
$ws.Range("J86").Value = 36969.332
$ws.Range("K86").Value = 11153860
$ws.Range("L86").Value = 36969.332
$ws.Range("M86").Value = -11152737
$ws.Range("N86").Value = -39215.332
$ws.Range("H87").Value = 19900
$ws.Range("J87").Value = 19900
$ws.Range("L87").Value = 19900
$ws.Range("N87").Value = -22272
$ws.Range("H89").Value = 4483726
$ws.Range("I89").Value = 11153860
$ws.Range("J89").Value = 36969.332
$ws.Range("K89").Value = 55769300
$ws.Range("L89").Value = 184846.66
$ws.Range("M89").Value = -55763684
$ws.Range("N89").Value = -196078.66
$ws.Range("H90").Value = 19900
$ws.Range("J90").Value = 19900
$ws.Range("L90").Value = 59700
$ws.Range("N90").Value = -71556
$ws.Range("H105").Value = 812.5
$ws.Range("I105").Value = 786.36365
$ws.Range("J105").Value = 1100
$ws.Range("K105").Value = 786.36365
$ws.Range("L105").Value = 1100
$ws.Range("M105").Value = 960.63635
$ws.Range("N105").Value = -4594
$ws.Range("H113").Value = 166668220
$ws.Range("I113").Value = 166668220
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 166668220
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -166666050
$ws.Range("H120").Value = 18799.334
$ws.Range("J120").Value = 18799.334
$ws.Range("L120").Value = 18799.334
$ws.Range("N120").Value = -26057.334

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2054.4546
$ws.Range("I34").Value = 250
$ws.Range("J34").Value = 2455.4443
$ws.Range("K34").Value = 750
$ws.Range("L34").Value = 7366.3329
$ws.Range("M34").Value = -666
$ws.Range("N34").Value = -7534.3329
$ws.Range("H99").Value = 1494.4166
$ws.Range("I99").Value = 760
$ws.Range("J99").Value = 2019
$ws.Range("K99").Value = 2280
$ws.Range("L99").Value = 6057
$ws.Range("M99").Value = -34
$ws.Range("N99").Value = -10549
$ws.Range("H131").Value = 10527363
$ws.Range("I131").Value = 125000370
$ws.Range("J131").Value = 1109.6093
$ws.Range("K131").Value = 375001110
$ws.Range("L131").Value = 3328.8279
$ws.Range("M131").Value = -374996070
$ws.Range("N131").Value = -13408.8279
$ws.Range("H132").Value = 954.5625
$ws.Range("I132").Value = 905.61536
$ws.Range("K132").Value = 8150.53824
$ws.Range("M132").Value = -5620.53824
$ws.Range("H133").Value = 4760
$ws.Range("I133").Value = 3412
$ws.Range("J133").Value = 5883.3335
$ws.Range("K133").Value = 10236
$ws.Range("L133").Value = 17650.0005
$ws.Range("M133").Value = -5176
$ws.Range("N133").Value = -27770.0005
$ws.Range("H139").Value = 1810.5454
$ws.Range("I139").Value = 1892.579
$ws.Range("J139").Value = 1699.2142
$ws.Range("K139").Value = 5677.737
$ws.Range("L139").Value = 5097.642599999999
$ws.Range("M139").Value = -537.7370000000001
$ws.Range("N139").Value = -15377.6426

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 16000
$ws.Range("J94").Value = 16000
$ws.Range("L94").Value = 16000
$ws.Range("N94").Value = -17352

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 14000
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 15333.333
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 15333.333
$ws.Range("M42").Value = -9437
$ws.Range("N42").Value = -16459.333
$ws.Range("H49").Value = 14000
$ws.Range("I49").Value = 10000
$ws.Range("J49").Value = 15333.333
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 15333.333
$ws.Range("M49").Value = -9853
$ws.Range("N49").Value = -15627.333
$ws.Range("H98").Value = 22000
$ws.Range("J98").Value = 22000
$ws.Range("L98").Value = 22000
$ws.Range("N98").Value = -27990
$ws.Range("H132").Value = 2410.257
$ws.Range("I132").Value = 2090.0417
$ws.Range("K132").Value = 6270.125100000001
$ws.Range("M132").Value = -3740.125100000001
$ws.Range("H133").Value = 46845.2
$ws.Range("J133").Value = 46845.2
$ws.Range("L133").Value = 46845.2
$ws.Range("N133").Value = -51905.2

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 35831.668
$ws.Range("J123").Value = 35831.668
$ws.Range("L123").Value = 35831.668
$ws.Range("N123").Value = -45631.668
